$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.078.31'
$ws.Range("E2").Value = '  +0.39%  '
$ws.Range("D3").Value = '1.824.91'
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.39%  '
$ws.Range("D5").Value = '312.29'
$ws.Range("E5").Value = '  +0.43%  '
$ws.Range("E6").Value = '  +0.35%  '
$ws.Range("D7").Value = '0.4690'
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("D8").Value = '0.3649'
$ws.Range("E8").Value = '  -0.47%  '
$ws.Range("D9").Value = '0.07390'
$ws.Range("E9").Value = '  +0.47%  '
$ws.Range("D10").Value = '0.8782'
$ws.Range("E10").Value = '  +0.53%  '
$ws.Range("D11").Value = '20.23'
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("D12").Value = '1.874.71'
$ws.Range("E12").Value = '  +2.59%  '
$ws.Range("D13").Value = '0.07536'
$ws.Range("E13").Value = '  +5.11%  '
$ws.Range("D14").Value = '5.373'
$ws.Range("E14").Value = '  -0.87%  '
$ws.Range("D15").Value = '92.79'
$ws.Range("E15").Value = '  +1.28%  '
$ws.Range("D16").Value = '6.524'
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("D17").Value = '1.005'
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").Value = '0.000008719'
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("E19").Value = '  +0.40%  '
$ws.Range("D20").Value = '27.532.24'
$ws.Range("E20").Value = '  +2.02%  '
$ws.Range("D21").Value = '14.61'
$ws.Range("E21").Value = '  -0.43%  '
$ws.Range("D22").Value = '5.237'
$ws.Range("E22").Value = '  -1.00%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = '2.082.68'
$ws.Range("E24").Value = '  +1.82%  '
$ws.Range("D25").Value = '1.876'
$ws.Range("E25").Value = '  -0.79%  '
$ws.Range("D26").Value = '151.40'
$ws.Range("E26").Value = '  +0.38%  '
$ws.Range("D27").Value = '18.51'
$ws.Range("E27").Value = '  +0.61%  '
$ws.Range("D28").Value = '2.132'
$ws.Range("E28").Value = '  -0.53%  '
$ws.Range("D29").Value = '5.165'
$ws.Range("E29").Value = '  -1.30%  '
$ws.Range("D30").Value = '116.42'
$ws.Range("E30").Value = '  -0.16%  '
$ws.Range("D31").Value = '0.08910'
$ws.Range("E31").Value = '  +0.30%  '
$ws.Range("D32").Value = '0.7430'
$ws.Range("E32").Value = '  -1.49%  '
$ws.Range("D33").Value = '1.162'
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D34").Value = '4.507'
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("D35").Value = '2.943'
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  +0.36%  '
$ws.Range("D37").Value = '2.534'
$ws.Range("E37").Value = '  +6.82%  '
$ws.Range("D38").Value = '1.090'
$ws.Range("D39").Value = '0.05298'
$ws.Range("E39").Value = '  -0.28%  '
$ws.Range("D40").Value = '0.01931'
$ws.Range("E40").Value = '  -0.85%  '
$ws.Range("D41").Value = '7.318'
$ws.Range("E41").Value = '  +1.84%  '
$ws.Range("D42").Value = '2.931'
$ws.Range("E42").Value = '  -1.63%  '
$ws.Range("D43").Value = '0.5244'
$ws.Range("E43").Value = '  -0.97%  '
$ws.Range("E44").Value = '  -0.76%  '
$ws.Range("D45").Value = '8.363'
$ws.Range("E45").Value = '  -1.16%  '
$ws.Range("D46").Value = '0.4890'
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").Value = '10.43'
$ws.Range("E47").Value = '  -0.80%  '
$ws.Range("E48").Value = '  +0.37%  '
$ws.Range("D49").Value = '104.24'
$ws.Range("E49").Value = '  +1.15%  '
$ws.Range("D50").Value = '1.649'
$ws.Range("E50").Value = '  -0.78%  '
$ws.Range("D51").Value = '0.06268'
$ws.Range("E51").Value = '  -0.46%  '
